$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(130, 130, 80, 75, 743, 560, 7200, 10875)
    3 = @(80, 110, 75, 75, 443, 255, 6750, 10875)
    4 = @(80, 130, 60, 110, 603, 400, 5400, 15950)
    5 = @(90, 130, 80, 110, 723, 410, 7200, 15950)
    6 = @(100, 160, 60, 130, 343, 400, 5400, 18850)
    7 = @(130, 140, 150, 90, 500, 400, 10500, 13050)
    8 = @(80, 110, 150, 110, 500, 400, 3000, 15950)
    9 = @(90, 80, 150, 110, 500, 400, 4500, 7550)
    10 = @(70, 110, 180, 75, 500, 400, -300, 10875)
    11 = @(100, 95, 75, 170, 510, 400, 6750, 3650)
    12 = @(70, 110, 60, 130, 630, 400, 5400, 13250)
    13 = @(100, 160, 150, 110, 610, 400, 6000, 15950)
    14 = @(90, 95, 180, 75, 600, 400, 2700, 10875)
    15 = @(80, 95, 60, 75, 330, 400, 5400, 10875)
    16 = @(70, 95, 180, 170, 500, 515, -300, 3650)
    17 = @(80, 160, 80, 110, 620, 400, 7200, 15950)
    18 = @(130, 95, 75, 90, 790, 400, 6750, 13050)
    19 = @(100, 130, 60, 110, 450, 455, 5400, 15950)
    20 = @(120, 130, 80, 110, 570, 465, 7200, 15950)
    21 = @(100, 130, 60, 110, 730, 455, 5400, 15950)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 2  # column B is index 2
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}
